$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=1.02; "C"=1.048804318342366; "D"=1.058095825123549; "E"=1.05550019707496; "F"=1.064786602124767; "I"=1.027286536288927; "J"=1.053846126285822; "K"=1.06082899501607; "L"=1.058240484293667; "M"=1.067501598779448; "N"=1.055342708662988 }
    3 = @{ "B"=1.02; "C"=1.050456132646184; "D"=1.059717873932031; "E"=1.057022838121862; "F"=1.066546039456642; "I"=1.027488238018589; "J"=1.055143524599395; "K"=1.062262960728831; "L"=1.059574771125205; "M"=1.069073947467342; "N"=1.05664194943098 }
    4 = @{ "B"=1.02; "C"=1.051516129805429; "D"=1.06075716925474; "E"=1.057998493754565; "F"=1.067669527853223; "I"=1.027608963341159; "J"=1.055973896163518; "K"=1.063180096727854; "L"=1.060428070540326; "M"=1.070075956318342; "N"=1.057473500217969 }
    5 = @{ "B"=1.02; "C"=1.051959665975854; "D"=1.061191658368259; "E"=1.05840639033651; "F"=1.068138292474475; "I"=1.027657380204171; "J"=1.056320823238079; "K"=1.063563118958202; "L"=1.060784412506436; "M"=1.070493547951727; "N"=1.057820919968768 }
    6 = @{ "B"=1.02; "C"=1.052034016108672; "D"=1.061264469316834; "E"=1.058474745698517; "F"=1.068216793005345; "I"=1.027665372821272; "J"=1.056378947796693; "K"=1.063627281826349; "L"=1.060844104771733; "M"=1.070563450207883; "N"=1.057879127070923 }
    7 = @{ "B"=1.02; "C"=1.051522064518176; "D"=1.060762984430582; "E"=1.058003952967518; "F"=1.067675805408264; "I"=1.02760961945904; "J"=1.055978540278783; "K"=1.063185224639523; "L"=1.060432841336066; "M"=1.070081550505112; "N"=1.05747815092841 }
    8 = @{ "B"=1.02; "C"=1.049364404379915; "D"=1.058646150982085; "E"=1.056016784737586; "F"=1.065384338259888; "I"=1.027356734046798; "J"=1.054286494425345; "K"=1.061315850902722; "L"=1.058693514910569; "M"=1.068036193157926; "N"=1.05578370217578 }
    9 = @{ "B"=1.02; "C"=1.045493282540499; "D"=1.054835871126434; "E"=1.052440327654229; "F"=1.061229921903791; "I"=1.026835758496609; "J"=1.05123376560083; "K"=1.057938203828243; "L"=1.055550191432082; "M"=1.064312334914406; "N"=1.052726638126182 }
    10 = @{ "B"=1.02; "C"=1.042864102964517; "D"=1.052239693480478; "E"=1.050003765052644; "F"=1.058379291337064; "I"=1.026437160985646; "J"=1.049149036313858; "K"=1.055628255595466; "L"=1.053400064795216; "M"=1.061746791731907; "N"=1.0506389482845 }
    11 = @{ "B"=1.02; "C"=1.041713677545962; "D"=1.051101739859867; "E"=1.048935846335879; "F"=1.057125086563524; "I"=1.026252243156335; "J"=1.048234149395155; "K"=1.054613747017149; "L"=1.052455647816042; "M"=1.060615602918453; "N"=1.049722762121433 }
    12 = @{ "B"=1.02; "C"=1.041284516768045; "D"=1.050676936031618; "E"=1.048537197981402; "F"=1.056656177412334; "I"=1.026181690060154; "J"=1.047892451209734; "K"=1.054234724058887; "L"=1.052102796252268; "M"=1.060192325800343; "N"=1.049380578685396 }
    13 = @{ "B"=1.02; "C"=1.041376657252431; "D"=1.050768154450273; "E"=1.0486227995127; "F"=1.056756898612233; "I"=1.026196908651096; "J"=1.047965831807091; "K"=1.054316125554137; "L"=1.052178577681205; "M"=1.060283261514093; "N"=1.049454063491619 }
    14 = @{ "B"=1.02; "C"=1.04167824080316; "D"=1.051066668901698; "E"=1.048902934494595; "F"=1.057086388778319; "I"=1.026246449380346; "J"=1.048205942857021; "K"=1.054582461778609; "L"=1.052426523116504; "M"=1.060580678330698; "N"=1.04969451552678 }
    15 = @{ "B"=1.02; "C"=1.041863811049387; "D"=1.051250311573437; "E"=1.049075271799178; "F"=1.057288993728594; "I"=1.026276725274696; "J"=1.048353634453407; "K"=1.054746268891255; "L"=1.052579017200067; "M"=1.060763513759093; "N"=1.049842416862173 }
    16 = @{ "B"=1.02; "C"=1.042940197044367; "D"=1.052314921334672; "E"=1.050074364576433; "F"=1.058462105445785; "I"=1.026449172492655; "J"=1.049209494432673; "K"=1.055695280577685; "L"=1.053462456975878; "M"=1.06182143288441; "N"=1.050699492260779 }
    17 = @{ "B"=1.02; "C"=1.043612150383544; "D"=1.052978999287257; "E"=1.050697594050364; "F"=1.059192609861862; "I"=1.026554034724548; "J"=1.049743064324509; "K"=1.056286716094883; "L"=1.054013000640728; "M"=1.062479567228624; "N"=1.051233819883062 }
    18 = @{ "B"=1.02; "C"=1.044002937130818; "D"=1.053365017197873; "E"=1.051059874336846; "F"=1.059616788277543; "I"=1.026614011233293; "J"=1.050053113545975; "K"=1.056630315074434; "L"=1.054332833552768; "M"=1.062861491841754; "N"=1.051544309409986 }
    19 = @{ "B"=1.02; "C"=1.044135991133596; "D"=1.053496415495051; "E"=1.051183193654496; "F"=1.059761099469125; "I"=1.026634260629912; "J"=1.050158634461549; "K"=1.056747241436066; "L"=1.054441670709686; "M"=1.062991388443327; "N"=1.051649980177364 }
    20 = @{ "B"=1.02; "C"=1.043540175637059; "D"=1.052907887639962; "E"=1.050630855871737; "F"=1.059114431872305; "I"=1.026542906970546; "J"=1.049685938928336; "K"=1.056223403225172; "L"=1.053954066220066; "M"=1.062409158141655; "N"=1.051176613362274 }
    21 = @{ "B"=1.02; "C"=1.041589483152185; "D"=1.050978822653564; "E"=1.048820496613536; "F"=1.056989446576506; "I"=1.026231912518285; "J"=1.048135288027346; "K"=1.054504093189688; "L"=1.052353566371638; "M"=1.060493182607434; "N"=1.049623760359143 }
    22 = @{ "B"=1.02; "C"=1.04035232820555; "D"=1.049753669680425; "E"=1.047670798323042; "F"=1.05563575542647; "I"=1.026025571364273; "J"=1.047149502042212; "K"=1.053410405331989; "L"=1.051335369033051; "M"=1.059270549409009; "N"=1.048636574444839 }
    23 = @{ "B"=1.02; "C"=1.041009194263411; "D"=1.050404325434018; "E"=1.048281375358128; "F"=1.056355063617472; "I"=1.026135986386036; "J"=1.04767312547901; "K"=1.053991407663165; "L"=1.051876276550971; "M"=1.059920413963197; "N"=1.049160941486977 }
    24 = @{ "B"=1.02; "C"=1.043572701485409; "D"=1.052940024032333; "E"=1.050661015808243; "F"=1.059149763047318; "I"=1.026547938793402; "J"=1.049711755058927; "K"=1.056252015833505; "L"=1.053980700131814; "M"=1.062440979012555; "N"=1.051202466154732 }
    25 = @{ "B"=1.02; "C"=1.046502443619029; "D"=1.055830627671226; "E"=1.053373988506274; "F"=1.062318014848719; "I"=1.02697943277255; "J"=1.052031573115911; "K"=1.058821508967545; "L"=1.056372293002348; "M"=1.065289473043913; "N"=1.053525578619444 }
}

foreach ($rowKey in $data.Keys) {
    $rowNum = [int]$rowKey
    $rowVals = $data[$rowKey]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
